# Auto-generated Excel COM-interop edit script
# Applies the "Update countries & provincias Spain" data refresh:
#  - reorders several province/city labels (shared-string shuffle in source diff)
#  - updates the Casos totales/activos/Recuperados/Muertes counts
#  - updates the "Datos actualizados..." timestamp string in A1
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(1, 1).Value = "Datos actualizados a 18 de Abril de 2020 a las 12:23"
$ws.Cells.Item(4, 2).Value = 51993
$ws.Cells.Item(4, 3).Value = 29436
$ws.Cells.Item(4, 4).Value = 15550
$ws.Cells.Item(4, 5).Value = 7007
$ws.Cells.Item(5, 2).Value = 40988
$ws.Cells.Item(5, 3).Value = 19088
$ws.Cells.Item(5, 4).Value = 17777
$ws.Cells.Item(5, 5).Value = 4123
$ws.Cells.Item(6, 2).Value = 8013
$ws.Cells.Item(6, 3).Value = 1456
$ws.Cells.Item(6, 4).Value = 6142
$ws.Cells.Item(6, 5).Value = 415
$ws.Cells.Item(7, 2).Value = 6667
$ws.Cells.Item(7, 3).Value = 6144
$ws.Cells.Item(7, 4).Value = 4953
$ws.Cells.Item(7, 5).Value = 500
$ws.Cells.Item(8, 2).Value = 6116
$ws.Cells.Item(8, 3).Value = 3600
$ws.Cells.Item(8, 4).Value = 10545
$ws.Cells.Item(8, 5).Value = 708
$ws.Cells.Item(9, 2).Value = 5131
$ws.Cells.Item(9, 3).Value = 2065
$ws.Cells.Item(9, 4).Value = 2576
$ws.Cells.Item(9, 5).Value = 490
$ws.Cells.Item(10, 2).Value = 4433
$ws.Cells.Item(10, 3).Value = 892
$ws.Cells.Item(10, 4).Value = 3193
$ws.Cells.Item(10, 5).Value = 348
$ws.Cells.Item(11, 1).Value = "La Rioja"
$ws.Cells.Item(11, 2).Value = 4098
$ws.Cells.Item(11, 3).Value = 1729
$ws.Cells.Item(11, 4).Value = 2107
$ws.Cells.Item(11, 5).Value = 262
$ws.Cells.Item(12, 1).Value = "Toledo"
$ws.Cells.Item(12, 2).Value = 3751
$ws.Cells.Item(12, 3).Value = 3600
$ws.Cells.Item(12, 4).Value = 10545
$ws.Cells.Item(12, 5).Value = 472
$ws.Cells.Item(13, 1).Value = "Albacete"
$ws.Cells.Item(13, 2).Value = 3673
$ws.Cells.Item(13, 3).Value = 3600
$ws.Cells.Item(13, 4).Value = 10545
$ws.Cells.Item(13, 5).Value = 347
$ws.Cells.Item(14, 1).Value = "Alacant/Alicante"
$ws.Cells.Item(14, 2).Value = 3476
$ws.Cells.Item(14, 3).Value = 1547
$ws.Cells.Item(14, 4).Value = 1549
$ws.Cells.Item(14, 5).Value = 380
$ws.Cells.Item(15, 1).Value = "Zaragoza"
$ws.Cells.Item(15, 2).Value = 3360
$ws.Cells.Item(15, 3).Value = 845
$ws.Cells.Item(15, 4).Value = 2089
$ws.Cells.Item(15, 5).Value = 426
$ws.Cells.Item(16, 2).Value = 3156
$ws.Cells.Item(16, 3).Value = 6144
$ws.Cells.Item(16, 4).Value = 4953
$ws.Cells.Item(16, 5).Value = 304
$ws.Cells.Item(17, 2).Value = 3067
$ws.Cells.Item(17, 3).Value = 1019
$ws.Cells.Item(17, 4).Value = 1813
$ws.Cells.Item(17, 5).Value = 235
$ws.Cells.Item(19, 1).Value = "Salamanca"
$ws.Cells.Item(19, 2).Value = 2437
$ws.Cells.Item(19, 3).Value = 733
$ws.Cells.Item(19, 4).Value = 1431
$ws.Cells.Item(19, 5).Value = 273
$ws.Cells.Item(20, 1).Value = "Malaga"
$ws.Cells.Item(20, 2).Value = 2429
$ws.Cells.Item(20, 3).Value = 809
$ws.Cells.Item(20, 4).Value = 1406
$ws.Cells.Item(20, 5).Value = 214
$ws.Cells.Item(21, 1).Value = "Sevilla"
$ws.Cells.Item(21, 2).Value = 2278
$ws.Cells.Item(21, 3).Value = 416
$ws.Cells.Item(21, 4).Value = 1665
$ws.Cells.Item(21, 5).Value = 197
$ws.Cells.Item(22, 1).Value = "Asturias"
$ws.Cells.Item(22, 2).Value = 2272
$ws.Cells.Item(22, 3).Value = 575
$ws.Cells.Item(22, 4).Value = 1510
$ws.Cells.Item(22, 5).Value = 187
$ws.Cells.Item(23, 1).Value = "Gipuzkoa/Guipuzcoa"
$ws.Cells.Item(23, 2).Value = 2266
$ws.Cells.Item(23, 3).Value = 6144
$ws.Cells.Item(23, 4).Value = 4953
$ws.Cells.Item(23, 5).Value = 188
$ws.Cells.Item(24, 2).Value = 2252
$ws.Cells.Item(24, 3).Value = 622
$ws.Cells.Item(24, 4).Value = 1469
$ws.Cells.Item(24, 5).Value = 161
$ws.Cells.Item(25, 1).Value = "Leon"
$ws.Cells.Item(25, 2).Value = 2212
$ws.Cells.Item(25, 3).Value = 998
$ws.Cells.Item(25, 4).Value = 927
$ws.Cells.Item(25, 5).Value = 287
$ws.Cells.Item(26, 1).Value = "Caceres"
$ws.Cells.Item(26, 2).Value = 2067
$ws.Cells.Item(26, 3).Value = 321
$ws.Cells.Item(26, 4).Value = 1442
$ws.Cells.Item(26, 5).Value = 304
$ws.Cells.Item(27, 1).Value = "A Coruña"
$ws.Cells.Item(27, 2).Value = 1969
$ws.Cells.Item(27, 3).Value = 333
$ws.Cells.Item(27, 4).Value = 1788
$ws.Cells.Item(27, 5).Value = 67
$ws.Cells.Item(28, 1).Value = "Granada"
$ws.Cells.Item(28, 2).Value = 1956
$ws.Cells.Item(28, 3).Value = 513
$ws.Cells.Item(28, 4).Value = 1252
$ws.Cells.Item(28, 5).Value = 191
$ws.Cells.Item(29, 2).Value = 1884
$ws.Cells.Item(29, 3).Value = 510
$ws.Cells.Item(29, 4).Value = 1230
$ws.Cells.Item(29, 5).Value = 144
$ws.Cells.Item(30, 2).Value = 1647
$ws.Cells.Item(30, 3).Value = 638
$ws.Cells.Item(30, 4).Value = 894
$ws.Cells.Item(30, 5).Value = 115
$ws.Cells.Item(32, 2).Value = 1450
$ws.Cells.Item(32, 3).Value = 606
$ws.Cells.Item(32, 4).Value = 685
$ws.Cells.Item(32, 5).Value = 159
$ws.Cells.Item(33, 2).Value = 1367
$ws.Cells.Item(33, 3).Value = 487
$ws.Cells.Item(33, 4).Value = 802
$ws.Cells.Item(33, 5).Value = 78
$ws.Cells.Item(34, 2).Value = 1312
$ws.Cells.Item(34, 3).Value = 3600
$ws.Cells.Item(34, 4).Value = 10545
$ws.Cells.Item(34, 5).Value = 176
$ws.Cells.Item(35, 2).Value = 1259
$ws.Cells.Item(35, 3).Value = 317
$ws.Cells.Item(35, 4).Value = 872
$ws.Cells.Item(35, 5).Value = 70
$ws.Cells.Item(36, 2).Value = 1257
$ws.Cells.Item(36, 3).Value = 412
$ws.Cells.Item(36, 4).Value = 716
$ws.Cells.Item(36, 5).Value = 129
$ws.Cells.Item(37, 2).Value = 1233
$ws.Cells.Item(37, 3).Value = 252
$ws.Cells.Item(37, 4).Value = 852
$ws.Cells.Item(37, 5).Value = 129
$ws.Cells.Item(38, 2).Value = 1179
$ws.Cells.Item(38, 3).Value = 285
$ws.Cells.Item(38, 4).Value = 802
$ws.Cells.Item(38, 5).Value = 92
$ws.Cells.Item(39, 1).Value = "Cuenca"
$ws.Cells.Item(39, 2).Value = 1145
$ws.Cells.Item(39, 3).Value = 3600
$ws.Cells.Item(39, 4).Value = 10545
$ws.Cells.Item(39, 5).Value = 149
$ws.Cells.Item(40, 1).Value = "Cadiz"
$ws.Cells.Item(40, 2).Value = 1110
$ws.Cells.Item(40, 3).Value = 276
$ws.Cells.Item(40, 4).Value = 763
$ws.Cells.Item(40, 5).Value = 71
$ws.Cells.Item(41, 1).Value = "Avila"
$ws.Cells.Item(41, 2).Value = 1075
$ws.Cells.Item(41, 3).Value = 437
$ws.Cells.Item(41, 4).Value = 531
$ws.Cells.Item(41, 5).Value = 107
$ws.Cells.Item(42, 1).Value = "Badajoz"
$ws.Cells.Item(42, 2).Value = 952
$ws.Cells.Item(42, 3).Value = 386
$ws.Cells.Item(42, 4).Value = 498
$ws.Cells.Item(42, 5).Value = 68
$ws.Cells.Item(43, 1).Value = "Aragon"
$ws.Cells.Item(43, 2).Value = 907
$ws.Cells.Item(43, 3).Value = 29
$ws.Cells.Item(43, 4).Value = 838
$ws.Cells.Item(43, 5).Value = 40
$ws.Cells.Item(45, 2).Value = 659
$ws.Cells.Item(45, 3).Value = 206
$ws.Cells.Item(45, 4).Value = 399
$ws.Cells.Item(45, 5).Value = 54
$ws.Cells.Item(47, 1).Value = "Zamora"
$ws.Cells.Item(47, 2).Value = 572
$ws.Cells.Item(47, 3).Value = 197
$ws.Cells.Item(47, 4).Value = 314
$ws.Cells.Item(47, 5).Value = 61
$ws.Cells.Item(48, 1).Value = "Huesca"
$ws.Cells.Item(48, 2).Value = 544
$ws.Cells.Item(48, 3).Value = 125
$ws.Cells.Item(48, 4).Value = 343
$ws.Cells.Item(48, 5).Value = 76
$ws.Cells.Item(49, 1).Value = "Teruel"
$ws.Cells.Item(49, 2).Value = 527
$ws.Cells.Item(49, 3).Value = 133
$ws.Cells.Item(49, 4).Value = 335
$ws.Cells.Item(49, 5).Value = 59
$ws.Cells.Item(50, 2).Value = 495
$ws.Cells.Item(50, 3).Value = 235
$ws.Cells.Item(50, 4).Value = 228
$ws.Cells.Item(50, 5).Value = 32
$ws.Cells.Item(51, 2).Value = 436
$ws.Cells.Item(51, 3).Value = 118
$ws.Cells.Item(51, 4).Value = 279
$ws.Cells.Item(51, 5).Value = 39
$ws.Cells.Item(52, 2).Value = 352
$ws.Cells.Item(52, 3).Value = 89
$ws.Cells.Item(52, 4).Value = 234
$ws.Cells.Item(52, 5).Value = 29
$ws.Cells.Item(54, 1).Value = "Ceuta"
$ws.Cells.Item(54, 2).Value = 110
$ws.Cells.Item(54, 3).Value = 54
$ws.Cells.Item(54, 4).Value = 52
$ws.Cells.Item(54, 5).Value = 4
$ws.Cells.Item(55, 1).Value = "Melilla"
$ws.Cells.Item(55, 2).Value = 104
$ws.Cells.Item(55, 3).Value = 30
$ws.Cells.Item(55, 4).Value = 72
$ws.Cells.Item(55, 5).Value = 2
$ws.Cells.Item(56, 2).Value = 73
$ws.Cells.Item(56, 3).Value = 21
$ws.Cells.Item(56, 4).Value = 49
$ws.Cells.Item(56, 5).Value = 3
$ws.Cells.Item(57, 2).Value = 69
$ws.Cells.Item(57, 3).Value = 18
$ws.Cells.Item(57, 4).Value = 51
$ws.Cells.Item(57, 5).Value = 48
$ws.Cells.Item(59, 2).Value = 24
$ws.Cells.Item(59, 3).Value = 18
$ws.Cells.Item(59, 4).Value = 6
$ws.Cells.Item(62, 1).Value = "Arroyo de la Luz"
$ws.Cells.Item(62, 2).Value = 7
$ws.Cells.Item(62, 3).Value = 0
$ws.Cells.Item(62, 4).Value = 7
$ws.Cells.Item(63, 1).Value = "La Gomera"
$ws.Cells.Item(63, 3).Value = 5
$ws.Cells.Item(63, 4).Value = 2
$ws.Cells.Item(64, 2).Value = 1
$ws.Cells.Item(64, 3).Value = 1
$ws.Cells.Item(64, 4).Value = 0
